$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# ---------------------------------------------------------------------------
# 1. Update value of B12 (57 -> 62)
# ---------------------------------------------------------------------------
$ws.Range("B12").Value = 62

# ---------------------------------------------------------------------------
# 2. Copy formatting from the still-untouched source cells into the new
#    destinations BEFORE we overwrite / clear anything, so the clipboard
#    always has the correct source format available.
# ---------------------------------------------------------------------------

# Row 17 F/G/H get the same look as the other data cells in that row (s=1)
$ws.Range("C17").Copy()
$ws.Range("F17").PasteSpecial($xlPasteFormats)
$ws.Range("G17").PasteSpecial($xlPasteFormats)
$ws.Range("H17").PasteSpecial($xlPasteFormats)

# Row 18 D (bold label) takes the look of the old "Ukupno: " label cell (A19)
$ws.Range("A19").Copy()
$ws.Range("D18").PasteSpecial($xlPasteFormats)

# Row 18 E/F/G/H (red summary numbers) take the look of the old B19:E19 cells
$ws.Range("B19").Copy()
$ws.Range("E18").PasteSpecial($xlPasteFormats)
$ws.Range("C19").Copy()
$ws.Range("F18").PasteSpecial($xlPasteFormats)
$ws.Range("D19").Copy()
$ws.Range("G18").PasteSpecial($xlPasteFormats)
$ws.Range("E19").Copy()
$ws.Range("H18").PasteSpecial($xlPasteFormats)

# Row 19 D (bold label) takes the look of the old "Postotak: " label cell (A20)
$ws.Range("A20").Copy()
$ws.Range("D19").PasteSpecial($xlPasteFormats)

# Row 19 E ("100% (ukupno)" wrap text) takes the look of the old B20 cell
$ws.Range("B20").Copy()
$ws.Range("E19").PasteSpecial($xlPasteFormats)

# Row 19 F/G/H (percentage numbers) take the look of the old C20:E20 cells
$ws.Range("C20").Copy()
$ws.Range("F19").PasteSpecial($xlPasteFormats)
$ws.Range("D20").Copy()
$ws.Range("G19").PasteSpecial($xlPasteFormats)
$ws.Range("E20").Copy()
$ws.Range("H19").PasteSpecial($xlPasteFormats)

# New row 19 A/B (graphs.js / 29) look like the other plain data rows (s=1)
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial($xlPasteFormats)
$ws.Range("B18").Copy()
$ws.Range("B19").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------------
# 3. Now that formatting has been copied, remove the stale content that
#    moved or is no longer needed.
# ---------------------------------------------------------------------------
$ws.Range("D17").Clear()
$ws.Range("E17").Clear()
$ws.Range("C18").Clear()
$ws.Range("C19").Clear()

# ---------------------------------------------------------------------------
# 4. Write the new values / formulas into their destinations.
# ---------------------------------------------------------------------------

# Row 17: HTML / CSS / JS headers moved from row 18 to F17:H17
$ws.Range("F17").Value = "HTML"
$ws.Range("G17").Value = "CSS"
$ws.Range("H17").Value = "JS"

# Row 18: "Ukupno: " label and totals (now referencing the extra row 19)
$ws.Range("D18").Value = "Ukupno: "
$ws.Range("E18").Formula = "=SUM(B1:B19)"
$ws.Range("F18").Formula = "=SUM(B1,B2,B3,B4,B5,B6,B7,B9,B12,B14,B16)"
$ws.Range("G18").Formula = "=SUM(B8,B10,B13,B17,B15)"
$ws.Range("H18").Formula = "=SUM(B11,B18,B19)"

# Row 19: new data row (graphs.js / 29) plus "Postotak: " percentages
$ws.Range("A19").Value = "graphs.js"
$ws.Range("B19").Value = 29
$ws.Range("D19").Value = "Postotak: "
$ws.Range("E19").Value = "100% (ukupno)"
$ws.Range("F19").Formula = "=(F18/E18)*100"
$ws.Range("G19").Formula = "=(G18/E18)*100"
$ws.Range("H19").Formula = "=(H18/E18)*100"

# ---------------------------------------------------------------------------
# 5. The old row 20 (previous percentage row) is no longer needed now that
#    its content lives in row 19 - remove the now-duplicate row entirely.
# ---------------------------------------------------------------------------
$ws.Rows(20).Delete()

# ---------------------------------------------------------------------------
# 6. Update the active selection to match the author's saved cursor (F19).
# ---------------------------------------------------------------------------
$ws.Range("F19").Select()
